$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D5: end time for MAT141 row
$ws.Range("D5").Value = 45689.76860903935

# Update B8: duration for MILLONIARA row
$ws.Range("B8").Value = 9

# Update D8: end time for MILLONIARA row
$ws.Range("D8").Value = 45689.76876608796

# Add new row 9
$ws.Range("A9").Value = "aaaaaa"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 45689.76918370945
$ws.Range("D9").Value = 45689.76918370945

# Apply same number format as other data rows (style index 2, numFmt 165) to C9/D9
$ws.Range("C9").NumberFormat = $ws.Range("C8").NumberFormat
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat
